$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Row, $Col, $Text) {
    $cell = $ws.Cells.Item($Row, $Col)
    $escaped = $Text -replace '"', '""'
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null
}

Set-TextValue 2 4 '20.562.47'
Set-TextValue 2 5 '  +1.70%  '
Set-TextValue 3 4 '1.473.77'
Set-TextValue 3 5 '  +1.93%  '
Set-TextValue 4 4 '1.011'
Set-TextValue 4 5 '  -0.66%  '
Set-TextValue 5 4 '0.9500'
Set-TextValue 5 5 '  +6.03%  '
Set-TextValue 6 4 '278.84'
Set-TextValue 6 5 '  +0.52%  '
Set-TextValue 7 4 '0.3626'
Set-TextValue 7 5 '  -1.60%  '
Set-TextValue 8 4 '0.3053'
Set-TextValue 8 5 '  -2.75%  '
Set-TextValue 9 4 '39.50'
Set-TextValue 9 5 '  +1.12%  '
Set-TextValue 10 4 '1.054'
Set-TextValue 10 5 '  +2.93%  '
Set-TextValue 11 4 '0.06654'
Set-TextValue 11 5 '  +2.21%  '
Set-TextValue 12 4 '1.004'
Set-TextValue 12 5 '  -0.65%  '
Set-TextValue 13 4 '5.508'
Set-TextValue 13 5 '  +1.73%  '
Set-TextValue 14 4 '17.98'
Set-TextValue 14 5 '  +3.30%  '
Set-TextValue 15 4 '6.200'
Set-TextValue 15 5 '  +1.33%  '
Set-TextValue 16 4 '0.9506'
Set-TextValue 16 5 '  +5.46%  '
Set-TextValue 17 4 '0.00001028'
Set-TextValue 17 5 '  +1.04%  '
Set-TextValue 18 4 '1.473.06'
Set-TextValue 18 5 '  +1.32%  '
Set-TextValue 19 4 '0.05940'
Set-TextValue 19 5 '  +6.17%  '
Set-TextValue 20 4 '69.24'
Set-TextValue 20 5 '  +2.73%  '
Set-TextValue 21 4 '5.480'
Set-TextValue 21 5 '  +0.27%  '
Set-TextValue 22 5 '  -0.06%  '
Set-TextValue 23 4 '11.07'
Set-TextValue 23 5 '  +0.16%  '
Set-TextValue 24 4 '2.276'
Set-TextValue 24 5 '  +0.90%  '
Set-TextValue 25 4 '20.581.25'
Set-TextValue 25 5 '  +0.58%  '
Set-TextValue 26 4 '142.89'
Set-TextValue 26 5 '  +5.72%  '
Set-TextValue 27 4 '2.108'
Set-TextValue 27 5 '  -3.77%  '
Set-TextValue 28 4 '17.19'
Set-TextValue 28 5 '  +1.15%  '
Set-TextValue 29 4 '1.633.06'
Set-TextValue 29 5 '  +1.27%  '
Set-TextValue 30 4 '113.47'
Set-TextValue 30 5 '  +2.24%  '
Set-TextValue 31 4 '3.951'
Set-TextValue 31 5 '  +7.76%  '
Set-TextValue 32 4 '4.998'
Set-TextValue 32 5 '  +2.43%  '
Set-TextValue 33 4 '0.8055'
Set-TextValue 33 5 '  -0.08%  '
Set-TextValue 34 4 '0.07948'
Set-TextValue 34 5 '  +3.48%  '
Set-TextValue 35 4 '1.508'
Set-TextValue 35 5 '  +6.75%  '
Set-TextValue 36 4 '1.212'
Set-TextValue 36 5 '  +5.92%  '
Set-TextValue 37 4 '0.05846'
Set-TextValue 37 5 '  -2.13%  '
Set-TextValue 38 4 '4.712'
Set-TextValue 38 5 '  +0.05%  '
Set-TextValue 39 4 '0.02046'
Set-TextValue 39 5 '  +1.69%  '
Set-TextValue 40 4 '10.34'
Set-TextValue 40 5 '  +0.79%  '
Set-TextValue 41 4 '0.9515'
Set-TextValue 41 5 '  +4.16%  '
Set-TextValue 42 4 '0.1874'
Set-TextValue 42 5 '  +1.81%  '
Set-TextValue 43 4 '7.391'
Set-TextValue 43 5 '  +9.56%  '
Set-TextValue 44 4 '0.5289'
Set-TextValue 44 5 '  +0.77%  '
Set-TextValue 45 4 '3.531'
Set-TextValue 45 5 '  -0.36%  '
Set-TextValue 46 4 '12.17'
Set-TextValue 46 5 '  +1.69%  '
Set-TextValue 47 4 '117.71'
Set-TextValue 47 5 '  -2.03%  '
Set-TextValue 48 4 '0.5177'
Set-TextValue 48 5 '  +0.59%  '
Set-TextValue 49 4 '1.813'
Set-TextValue 49 5 '  +2.37%  '
Set-TextValue 50 4 '0.06467'
Set-TextValue 50 5 '  +2.24%  '
Set-TextValue 51 4 '0.9840'
Set-TextValue 51 5 '  -1.40%  '

$excel.CutCopyMode = 0
